$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add Belgium data" - new location row appended to the table (row 115)
$ws.Range("A115").Value = "Touiouse"
$ws.Range("B115").Value = 43.6
$ws.Range("C115").Value = 1.46
$ws.Range("D115").Value = 1
$ws.Range("E115").Value = 1

# Move the active selection to where the editor left off (D109)
$ws.Range("D109").Select() | Out-Null
